$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows before row 68 (shifts old row 68.. down by 3).
$ws.Range("68:70").Insert()

# The sheet had 10 trailing blank template rows; after the insert there are
# 13 (10 + 3 carried down). Remove 3 of them so the used range / row count
# stays the same as before the edit (net effect: 3 rows "moved" from the
# bottom block of blanks to just above the old row 68).
$ws.Range("94:96").Delete()

# Give the 3 new rows the same formatting as the existing blank template
# rows (copy cell formats for columns A:L from one of the still-blank rows).
$ws.Range("A91:L91").Copy()
$ws.Range("A68:L70").PasteSpecial(-4122)

# Match row heights used by the template rows.
$ws.Rows.Item(68).RowHeight = 25.15
$ws.Rows.Item(69).RowHeight = 25.15
$ws.Rows.Item(70).RowHeight = 25.15

# Fill in the new feature rows.
$ws.Cells.Item(68,7).Value = "投资人员风控信息"
$ws.Cells.Item(69,7).Value = "分析决策表"
$ws.Cells.Item(70,7).Value = "个人首页"

# The AutoFilter range defined name needs to grow by the same 3 rows.
$names = $wb.Names
for ($i=1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=开发计划!`$A`$2:`$L`$73"
    }
}

# Reflect the user's new viewport/selection on the sheet.
$excel.ActiveWindow.ScrollRow = 64
$ws.Range("G68").Select()
